$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7228.875
$ws.Range("I19").Value = 9401.091
$ws.Range("K19").Value = 9401.091
$ws.Range("M19").Value = -9226.091

$ws.Range("H33").Value = 586.13635
$ws.Range("I33").Value = 710.1429000000001
$ws.Range("J33").Value = 369.125
$ws.Range("K33").Value = 710.1429000000001
$ws.Range("L33").Value = 369.125
$ws.Range("M33").Value = -481.1429000000001
$ws.Range("N33").Value = -827.125

$ws.Range("H100").Value = 5635.107
$ws.Range("I100").Value = 1765.8334
$ws.Range("J100").Value = 6690.364
$ws.Range("K100").Value = 1765.8334
$ws.Range("L100").Value = 6690.364
$ws.Range("M100").Value = -1224.8334
$ws.Range("N100").Value = -7772.364

$ws.Range("H125").Value = 3000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 27000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -31920

$ws.Range("H129").Value = 1750
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1750
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 5250
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -15250

$ws.Range("H137").Value = 4491.478
$ws.Range("I137").Value = 4812.75
$ws.Range("J137").Value = 3757.1428
$ws.Range("K137").Value = 14438.25
$ws.Range("L137").Value = 11271.4284
$ws.Range("M137").Value = -11888.25
$ws.Range("N137").Value = -16371.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 6000
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -1480
$ws.Range("N39").Value = -11040

$ws.Range("H45").Value = 849.8333
$ws.Range("I45").Value = 809.6
$ws.Range("K45").Value = 809.6
$ws.Range("M45").Value = -432.6

$ws.Range("H74").Value = 217543.2
$ws.Range("I74").Value = 271587.53
$ws.Range("J74").Value = 74711.71000000001
$ws.Range("K74").Value = 271587.53
$ws.Range("L74").Value = 74711.71000000001
$ws.Range("M74").Value = -270713.53
$ws.Range("N74").Value = -76459.71000000001

$ws.Range("H77").Value = 217543.2
$ws.Range("I77").Value = 271587.53
$ws.Range("J77").Value = 74711.71000000001
$ws.Range("K77").Value = 1357937.65
$ws.Range("L77").Value = 373558.55
$ws.Range("M77").Value = -1353569.65
$ws.Range("N77").Value = -382294.55

$ws.Range("H97").Value = 217.95238
$ws.Range("I97").Value = 229.77777
$ws.Range("J97").Value = 147
$ws.Range("K97").Value = 229.77777
$ws.Range("L97").Value = 147
$ws.Range("M97").Value = 266.22223
$ws.Range("N97").Value = -1139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7188.778
$ws.Range("I99").Value = 7837.375
$ws.Range("K99").Value = 7837.375
$ws.Range("M99").Value = -6339.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4752

$ws.Range("H33").Value = 2466.6667
$ws.Range("I33").Value = 2750
$ws.Range("J33").Value = 1900
$ws.Range("K33").Value = 2750
$ws.Range("L33").Value = 1900
$ws.Range("M33").Value = -2371
$ws.Range("N33").Value = -2658

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 847.05884
$ws.Range("I5").Value = 560
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1680
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -1568
$ws.Range("N5").Value = -9224

$ws.Range("H98").Value = 5260.3076
$ws.Range("J98").Value = 7487.1113
$ws.Range("L98").Value = 22461.3339
$ws.Range("N98").Value = -25457.3339

$ws.Range("H104").Value = 2795.923
$ws.Range("J104").Value = 2912
$ws.Range("L104").Value = 8736
$ws.Range("N104").Value = -13978

$ws.Range("H122").Value = 772.0303
$ws.Range("I122").Value = 381.6
$ws.Range("J122").Value = 1097.3889
$ws.Range("K122").Value = 3434.4
$ws.Range("L122").Value = 9876.500099999999
$ws.Range("M122").Value = -984.4000000000001
$ws.Range("N122").Value = -14776.5001

$ws.Range("H131").Value = 1482.9836
$ws.Range("I131").Value = 1104.6154
$ws.Range("J131").Value = 1585.4584
$ws.Range("K131").Value = 3313.8462
$ws.Range("L131").Value = 4756.3752
$ws.Range("M131").Value = 1726.1538
$ws.Range("N131").Value = -14836.3752

$ws.Range("H135").Value = 847.05884
$ws.Range("I135").Value = 560
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 5040
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -2505
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 15500.223
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 17375.25
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 17375.25
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = -17599.25

$ws.Range("H70").Value = 4250.926
$ws.Range("I70").Value = 4063.7932
$ws.Range("J70").Value = 4468
$ws.Range("K70").Value = 4063.7932
$ws.Range("L70").Value = 4468
$ws.Range("M70").Value = -3793.7932
$ws.Range("N70").Value = -5008

$ws.Range("H73").Value = 4250.926
$ws.Range("I73").Value = 4063.7932
$ws.Range("J73").Value = 4468
$ws.Range("K73").Value = 4063.7932
$ws.Range("L73").Value = 4468
$ws.Range("M73").Value = -3127.7932
$ws.Range("N73").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1177.591
$ws.Range("I16").Value = 1018.5128
$ws.Range("J16").Value = 2418.4
$ws.Range("K16").Value = 1018.5128
$ws.Range("L16").Value = 2418.4
$ws.Range("M16").Value = -848.5128
$ws.Range("N16").Value = -2758.4

$ws.Range("H132").Value = 11243.88
$ws.Range("I132").Value = 4124.875
$ws.Range("J132").Value = 14594
$ws.Range("K132").Value = 12374.625
$ws.Range("L132").Value = 43782
$ws.Range("M132").Value = -9844.625
$ws.Range("N132").Value = -48842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1080.4482
$ws.Range("I126").Value = 602.5217
$ws.Range("J126").Value = 2912.5
$ws.Range("K126").Value = 1807.5651
$ws.Range("L126").Value = 8737.5
$ws.Range("M126").Value = 662.4349
$ws.Range("N126").Value = -13677.5

$ws.Range("H132").Value = 2572.96
$ws.Range("I132").Value = 1161.7
$ws.Range("J132").Value = 3513.8
$ws.Range("K132").Value = 3485.1
$ws.Range("L132").Value = 10541.4
$ws.Range("M132").Value = -955.1000000000004
$ws.Range("N132").Value = -15601.4
